$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("sigma_010")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 27.45103973375281
$ws.Range("C2").Value = 28.99995409239356
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 27.47001208672447
$ws.Range("C3").Value = 28.98120564569585
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 27.4514745931836
$ws.Range("C4").Value = 28.98390008507995
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 27.45987433437825
$ws.Range("C5").Value = 28.98435864161655
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 27.46209955133791
$ws.Range("C6").Value = 28.96842979515579
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 27.45186209844839
$ws.Range("C7").Value = 28.98064374799339
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 27.45331013859491
$ws.Range("C8").Value = 28.97810287684981
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 27.46968888338561
$ws.Range("C9").Value = 28.98981616882067
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 27.46596940576997
$ws.Range("C10").Value = 28.97574312175779
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 27.46402022489861
$ws.Range("C11").Value = 28.98210435170081
$ws.Range("B12").Value = 27.45993510504745
$ws.Range("C12").Value = 28.98242585270642

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 19.64866192695948
$ws.Range("C2").Value = 25.56633709262722
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 19.64325667498629
$ws.Range("C3").Value = 25.53873280668903
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 19.65557126468072
$ws.Range("C4").Value = 25.6026981461732
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 19.65489040634082
$ws.Range("C5").Value = 25.54922003685539
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 19.64757837367189
$ws.Range("C6").Value = 25.53499126319608
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 19.64433111774049
$ws.Range("C7").Value = 25.56644980792897
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 19.65509330317094
$ws.Range("C8").Value = 25.53600006479947
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 19.65480343911532
$ws.Range("C9").Value = 25.54161828237438
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 19.64498065119335
$ws.Range("C10").Value = 25.51722175269063
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 19.64929363825734
$ws.Range("C11").Value = 25.55376956803781
$ws.Range("B12").Value = 19.64984607961166
$ws.Range("C12").Value = 25.55070388213722

$ws = $wb.Worksheets.Item("sigma_050")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 14.47904939829139
$ws.Range("C2").Value = 20.92379087772777
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 14.49210127248127
$ws.Range("C3").Value = 20.89577573656967
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 14.46984950534335
$ws.Range("C4").Value = 20.89497555576769
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 14.49322492543467
$ws.Range("C5").Value = 20.92430316553884
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 14.47744974752012
$ws.Range("C6").Value = 20.94186143709757
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 14.46620166359357
$ws.Range("C7").Value = 20.88580341287859
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 14.4569175079183
$ws.Range("C8").Value = 20.88379751100841
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 14.46458970736431
$ws.Range("C9").Value = 20.86756009837554
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 14.47102196411167
$ws.Range("C10").Value = 20.94703556754425
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 14.46574477025677
$ws.Range("C11").Value = 20.88756734196802
$ws.Range("B12").Value = 14.47361504623154
$ws.Range("C12").Value = 20.90524707044764
